# Update "Generate Report for Handback" timestamps on the zh-cn and de-de
# handback status sheets. The "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) values on row 2 of each sheet are
# refreshed with new report-generation timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-14 06:48:28"
$wsZhCn.Range("H2").Value = "2016-03-14 06:48:45"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-14 06:48:31"
$wsDeDe.Range("H2").Value = "2016-03-14 06:48:51"
